$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table 2 stats were recalculated for the obesity outcome (follow-up cap
# is now applied separately per outcome). Only the "Obese" column counts
# (D/E for rows 4-5) change; all other labels/cells stay the same.
$ws.Range("D4").Value = "108,259 (99)"
$ws.Range("E4").Value = "  1,208 (1)"
$ws.Range("D5").Value = " 43,812 (96)"
$ws.Range("E5").Value = "  1,920 (4)"
